# Expand the political-party abbreviations used throughout the workbook
# so that cells showing "MLP", "PN " (and the combined "MLP+PN ") display
# their full descriptive names instead of the bare abbreviation.
#
#   "MLP"     -> "MLP - Malta Labour Party (Partit Laburista, MLP)"
#   "PN "     -> "PN - Nationalist Party (Partit Nazzjonalista, PN)"
#   "MLP+PN " -> "<expanded MLP>+<expanded PN>"
#
# NOTE: reading with `.Value` in this runtime yields an unusable wrapper
# object when interpolated/compared, so `.Value2` is used for comparisons;
# `.Value` is used (safely) for writing the new text back.

$wb = $excel.ActiveWorkbook

$mlpOld = "MLP"
$mlpNew = "MLP - Malta Labour Party (Partit Laburista, MLP)"

$pnOld = "PN "
$pnNew = "PN - Nationalist Party (Partit Nazzjonalista, PN)"

$comboOld = "MLP+PN "
$comboNew = "$mlpNew+$pnNew"

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $rowCount = $used.Rows.Count
    $colCount = $used.Columns.Count

    for ($r = 1; $r -le $rowCount; $r++) {
        for ($c = 1; $c -le $colCount; $c++) {
            $cell = $used.Cells.Item($r, $c)
            $val = $cell.Value2

            if ($val -eq $comboOld) {
                $cell.Value = $comboNew
            } elseif ($val -eq $mlpOld) {
                $cell.Value = $mlpNew
            } elseif ($val -eq $pnOld) {
                $cell.Value = $pnNew
            }
        }
    }
}
